$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the role/description text in column C down to just the job title
$ws.Range("C1").Value = "CEO"
$ws.Range("C2").Value = "Director of Marketing"
$ws.Range("C3").Value = "Director of Development"
$ws.Range("C4").Value = "Content Manager"
$ws.Range("C5").Value = "Analytics Specialist"
$ws.Range("C6").Value = "Creative Director"
$ws.Range("C7").Value = "Project Manager"

# Shrink column C width now that the long descriptions are gone
$ws.Columns.Item(3).ColumnWidth = 22.7
